# Re-run of the backward-elimination export: the statsmodels OLS summary
# text embedded in column B row 2 of every step's sheet gets its
# "Date:"/"Time:" stamp refreshed to match the new export run
# (Thu 02 Jan 2020 -> Sun 05 Jan 2020, ~20:48 -> ~21:22).
$wb = $excel.ActiveWorkbook

$newDate = "Sun, 05 Jan 2020"
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value()

    if ($null -eq $text) {
        continue
    }

    if ($i -eq 1) {
        $newTime = "21:22:36"
    } else {
        $newTime = "21:22:37"
    }

    $text = $text.Replace("Thu, 02 Jan 2020", $newDate)
    $text = $text.Replace("20:48:57", $newTime)
    $text = $text.Replace("20:48:58", $newTime)

    $row = $ws.Rows.Item(2)
    $origRowHeight = $row.RowHeight()

    $cell.Value = $text

    # Re-assigning a wrapped cell's text re-triggers that row's auto-fit
    # (the cell style has wrapText=1); restore the original, already
    # maxed-out row height so the row geometry doesn't drift.
    $row.RowHeight = $origRowHeight
}
